$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Sheet "Normalize to 3NF" (sheet2): blank out the worked sample data,
# leaving only the styled, empty template cells behind.
# ---------------------------------------------------------------------
$ws2.Range("A2").ClearContents()
$ws2.Range("A3:E3").ClearContents()
$ws2.Range("A4").ClearContents()
$ws2.Range("B4").ClearContents()
$ws2.Range("C4").ClearContents()
$ws2.Range("E4").ClearContents()
$ws2.Range("A5").ClearContents()
$ws2.Range("B5").ClearContents()
$ws2.Range("C5").ClearContents()
$ws2.Range("E5").ClearContents()
$ws2.Range("A6").ClearContents()
$ws2.Range("B6").ClearContents()
$ws2.Range("C6").ClearContents()
$ws2.Range("E6").ClearContents()

# ---------------------------------------------------------------------
# Sheet "3NF Solution" (sheet3): fix a typo, then add the Orders /
# Products / OrderProducts normalized tables below the existing
# Customers table.
# ---------------------------------------------------------------------
$ws3.Range("B9").Value = "Age column removed as it was transitively dependent on the primary key through the Data of Birth column"

# Orders table heading (bold "Orders" + plain " table")
$ws3.Range("A12").Value = "Orders table"
$ws3.Range("A12").Characters(1, 6).Font.Bold = $true
$ws3.Range("A12").Characters(7, 6).Font.Bold = $false

# Products table heading (cell is bold by default; "Products" run plain, " table" run explicit)
$ws3.Range("E12").Value = "Products table"
$ws3.Range("E12").Font.Bold = $true
$ws3.Range("E12").Characters(10, 6).Font.Bold = $false

# Header rows, copying the existing header style (bold + shaded fill)
$ws3.Range("A3").Copy()
$ws3.Range("A13:C13").PasteSpecial(-4122)
$ws3.Range("E13:F13").PasteSpecial(-4122)
$ws3.Range("A20:C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("A13").Value = "OrderID"
$ws3.Range("B13").Value = "OrderDate"
$ws3.Range("C13").Value = "CustomerID"
$ws3.Range("E13").Value = "ProductID"
$ws3.Range("F13").Value = "Price"

# Orders data rows
$ws3.Range("A14").Value = 1
$ws3.Range("B14").Value = 42371
$ws3.Range("B14").NumberFormat = "[$-409]dd\-mmm\-yy;@"
$ws3.Range("C14").Value = 101
$ws3.Range("E14").Value = "CB-2903"
$ws3.Range("F14").Value = 12.99
$ws3.Range("F14").NumberFormat = "0.00"

$ws3.Range("A15").Value = 2
$ws3.Range("B15").Value = 42371
$ws3.Range("B15").NumberFormat = "[$-409]dd\-mmm\-yy;@"
$ws3.Range("C15").Value = 163
$ws3.Range("E15").Value = "BA-3827"
$ws3.Range("F15").Value = 1.5
$ws3.Range("F15").NumberFormat = "0.00"

$ws3.Range("A16").Value = 3
$ws3.Range("B16").Value = 42372
$ws3.Range("B16").NumberFormat = "[$-409]dd\-mmm\-yy;@"
$ws3.Range("C16").Value = 302
$ws3.Range("E16").Value = "BA-2349"
$ws3.Range("F16").Value = 5.99
$ws3.Range("F16").NumberFormat = "0.00"

$ws3.Range("E17").Value = "BA-2903"
$ws3.Range("F17").Value = 10
$ws3.Range("F17").NumberFormat = "0.00"

# OrderProducts table heading (cell bold by default; "OrderProducts" run plain, " table" run explicit)
$ws3.Range("A19").Value = "OrderProducts table"
$ws3.Range("A19").Font.Bold = $true
$ws3.Range("A19").Characters(14, 6).Font.Bold = $false

$ws3.Range("A20").Value = "OrderID"
$ws3.Range("B20").Value = "ProductID"
$ws3.Range("C20").Value = "Quantity"

$ws3.Range("A21").Value = 1
$ws3.Range("B21").Value = "CB-2903"
$ws3.Range("C21").Value = 1

$ws3.Range("A22").Value = 1
$ws3.Range("B22").Value = "BA-3827"
$ws3.Range("C22").Value = 2

$ws3.Range("A23").Value = 2
$ws3.Range("B23").Value = "BA-3827"
$ws3.Range("C23").Value = 1

$ws3.Range("A24").Value = 2
$ws3.Range("B24").Value = "BA-2349"
$ws3.Range("C24").Value = 1

$ws3.Range("A25").Value = 2
$ws3.Range("B25").Value = "BA-2908"
$ws3.Range("C25").Value = 1

$ws3.Range("A26").Value = 3
$ws3.Range("B26").Value = "BE-2349"
$ws3.Range("C26").Value = 1

# ---------------------------------------------------------------------
# Selections: update per-sheet cursor position. Activate sheet1, then
# sheet2, then sheet3 last so sheet3 ends up the active tab (as before).
# ---------------------------------------------------------------------
$ws1.Range("A1:B4").Select()
$ws2.Range("C16").Select()
$ws3.Range("I8").Select()
